$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.962.89"
$ws.Range("E2").Value = "  +4.76%  "

$ws.Range("D3").Value = "2.606.42"
$ws.Range("E3").Value = "  +5.34%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.24%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +1.90%  "

$ws.Range("D9").Value = "2.605.76"
$ws.Range("E9").Value = "  +5.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +16.24%  "

$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("E12").Value = "  +4.61%  "

$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.98%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.039.25"
$ws.Range("E15").Value = "  +3.83%  "

$ws.Range("E16").Value = "  +8.15%  "

$ws.Range("D17").Value = "70.991.60"
$ws.Range("E17").Value = "  +4.82%  "

$ws.Range("D18").Value = "2.595.08"
$ws.Range("E18").Value = "  +4.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "376.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.17%  "

$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("E26").Value = "  +12.35%  "

$ws.Range("E27").Value = "  +10.20%  "

$ws.Range("D28").Value = "2.741.20"
$ws.Range("E28").Value = "  +6.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").Value = "0.0₃0950"
$ws.Range("E30").Value = "  +6.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "530.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.16%  "

$ws.Range("E33").Value = "  +7.06%  "

$ws.Range("E34").Value = "  +4.67%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.58%  "

$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.61%  "

$ws.Range("E40").Value = "  +6.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.11%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").Value = "  +6.36%  "

$ws.Range("E45").Value = "  +1.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.75%  "

$ws.Range("E48").Value = "  +4.41%  "

$ws.Range("D49").Value = "0.0₆0268"
$ws.Range("E49").Value = "  +6.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.532"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.65%  "

$ws.Range("E51").Value = "  +7.28%  "

